$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 132
$ws.Range("H132").Value = 2030.1034
$ws.Range("J132").Value = 632
$ws.Range("L132").Value = 1896
$ws.Range("N132").Value = -6956
# row 137
$ws.Range("H137").Value = 1893.1096
$ws.Range("I137").Value = 1750
$ws.Range("K137").Value = 5250
$ws.Range("M137").Value = -2700
# row 138
$ws.Range("H138").Value = 114688.164
$ws.Range("I138").Value = 1206.5
$ws.Range("J138").Value = 128873.375
$ws.Range("K138").Value = 3619.5
$ws.Range("L138").Value = 386620.125
$ws.Range("M138").Value = 1520.5
$ws.Range("N138").Value = -396900.125

$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 3628.12
$ws.Range("I32").Value = 2902.2874
$ws.Range("K32").Value = 2902.2874
$ws.Range("M32").Value = -2615.2874
# row 52
$ws.Range("H52").Value = 94780
$ws.Range("J52").Value = 94780
$ws.Range("L52").Value = 94780
$ws.Range("N52").Value = -95416
# row 61
$ws.Range("H61").Value = 7220.868
$ws.Range("I61").Value = 3966.5435
$ws.Range("K61").Value = 3966.5435
$ws.Range("M61").Value = -3754.5435
# row 64
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
# row 67
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
# row 74
$ws.Range("H74").Value = 3505.1304
$ws.Range("I74").Value = 1434.7333
$ws.Range("J74").Value = 7387.125
$ws.Range("K74").Value = 1434.7333
$ws.Range("L74").Value = 7387.125
$ws.Range("M74").Value = -560.7333000000001
$ws.Range("N74").Value = -9135.125
# row 77
$ws.Range("H77").Value = 3505.1304
$ws.Range("I77").Value = 1434.7333
$ws.Range("J77").Value = 7387.125
$ws.Range("K77").Value = 7173.6665
$ws.Range("L77").Value = 36935.625
$ws.Range("M77").Value = -2805.6665
$ws.Range("N77").Value = -45671.625
# row 102
$ws.Range("H102").Value = 4652850.5
$ws.Range("I102").Value = 1770.2195
$ws.Range("K102").Value = 1770.2195
$ws.Range("M102").Value = -148.2194999999999
# row 132
$ws.Range("H132").Value = 1813.3334
$ws.Range("I132").Value = 1815.0952
$ws.Range("K132").Value = 5445.2856
$ws.Range("M132").Value = -2915.2856
# row 136
$ws.Range("H136").Value = 7220.868
$ws.Range("I136").Value = 3966.5435
$ws.Range("K136").Value = 11899.6305
$ws.Range("M136").Value = -9349.630500000001

$ws = $wb.Worksheets.Item("BSM")
# row 45
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()
# row 86
$ws.Range("H86").Value = 5553.375
$ws.Range("I86").Value = 4106
$ws.Range("K86").Value = 4106
$ws.Range("M86").Value = -2983
# row 89
$ws.Range("H89").Value = 5553.375
$ws.Range("I89").Value = 4106
$ws.Range("K89").Value = 20530
$ws.Range("M89").Value = -14914
# row 105
$ws.Range("H105").Value = 5061.3125
$ws.Range("I105").Value = 4375.5
$ws.Range("K105").Value = 4375.5
$ws.Range("M105").Value = -2628.5
# row 134
$ws.Range("H134").Value = 3607.1082
$ws.Range("I134").Value = 3369.4
$ws.Range("J134").Value = 7767
$ws.Range("K134").Value = 10108.2
$ws.Range("L134").Value = 23301
$ws.Range("M134").Value = -7573.200000000001
$ws.Range("N134").Value = -28371

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 2485.262
$ws.Range("J31").Value = 3882.6428
$ws.Range("L31").Value = 3882.6428
$ws.Range("N31").Value = -4472.6428
# row 34
$ws.Range("H34").Value = 2485.262
$ws.Range("J34").Value = 3882.6428
$ws.Range("L34").Value = 3882.6428
$ws.Range("N34").Value = -4286.6428
# row 57
$ws.Range("H57").Value = 19999
$ws.Range("J57").Value = 19999
$ws.Range("L57").Value = 19999
$ws.Range("N57").Value = -21119
# row 58
$ws.Range("H58").Value = 2850.423
$ws.Range("I58").Value = 2825.3845
$ws.Range("J58").Value = 2875.4614
$ws.Range("K58").Value = 2825.3845
$ws.Range("L58").Value = 2875.4614
$ws.Range("M58").Value = -2622.3845
$ws.Range("N58").Value = -3281.4614
# row 134
$ws.Range("H134").Value = 2702.6553
$ws.Range("I134").Value = 1494.591
$ws.Range("K134").Value = 4483.772999999999
$ws.Range("M134").Value = -1948.772999999999
# row 136
$ws.Range("H136").Value = 2850.423
$ws.Range("I136").Value = 2825.3845
$ws.Range("J136").Value = 2875.4614
$ws.Range("K136").Value = 8476.1535
$ws.Range("L136").Value = 8626.3842
$ws.Range("M136").Value = -5926.1535
$ws.Range("N136").Value = -13726.3842
# row 141
$ws.Range("H141").Value = 539954.8
$ws.Range("J141").Value = 539954.8
$ws.Range("L141").Value = 539954.8
$ws.Range("N141").Value = -550314.8

$ws = $wb.Worksheets.Item("CUL")
# row 2
$ws.Range("H2").Value = 99.95238000000001
$ws.Range("I2").Value = 72.454544
$ws.Range("K2").Value = 434.727264
$ws.Range("M2").Value = -321.727264
# row 51
$ws.Range("H51").Value = 2853.5
$ws.Range("I51").Value = 2853.5
$ws.Range("K51").Value = 8560.5
$ws.Range("M51").Value = -8100.5
# row 55
$ws.Range("H55").Value = 12343.625
$ws.Range("J55").Value = 13066.6
$ws.Range("L55").Value = 39199.8
$ws.Range("N55").Value = -39553.8
# row 56
$ws.Range("H56").Value = 9999.143
$ws.Range("I56").Value = 9999.143
$ws.Range("K56").Value = 9999.143
$ws.Range("M56").Value = -9469.143
# row 60
$ws.Range("H60").Value = 1954.3334
$ws.Range("I60").Value = 300.5
$ws.Range("J60").Value = 3608.1667
$ws.Range("K60").Value = 901.5
$ws.Range("L60").Value = 10824.5001
$ws.Range("M60").Value = -650.5
$ws.Range("N60").Value = -11326.5001
# row 63
$ws.Range("H63").Value = 999999
$ws.Range("I63").Value = 999999
$ws.Range("K63").Value = 2999997
$ws.Range("M63").Value = -2999248
# row 66
$ws.Range("H66").Value = 999999
$ws.Range("I66").Value = 999999
$ws.Range("K66").Value = 8999991
$ws.Range("M66").Value = -8996247
# row 113
$ws.Range("H113").Value = 2016.375
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2016.375
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 6049.125
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -10389.125
# row 121
$ws.Range("H121").Value = 3001.8462
$ws.Range("I121").Value = 1156.2727
$ws.Range("J121").Value = 13152.5
$ws.Range("K121").Value = 3468.8181
$ws.Range("L121").Value = 39457.5
$ws.Range("M121").Value = -2158.8181
$ws.Range("N121").Value = -42077.5

$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 5493.357
$ws.Range("I70").Value = 4888.778
$ws.Range("J70").Value = 6581.6
$ws.Range("K70").Value = 4888.778
$ws.Range("L70").Value = 6581.6
$ws.Range("M70").Value = -4618.778
$ws.Range("N70").Value = -7121.6
# row 73
$ws.Range("H73").Value = 5493.357
$ws.Range("I73").Value = 4888.778
$ws.Range("J73").Value = 6581.6
$ws.Range("K73").Value = 4888.778
$ws.Range("L73").Value = 6581.6
$ws.Range("M73").Value = -3952.778
$ws.Range("N73").Value = -8453.6
# row 80
$ws.Range("H80").Value = 3434.4583
$ws.Range("I80").Value = 2954
$ws.Range("K80").Value = 2954
$ws.Range("M80").Value = -1956
# row 83
$ws.Range("H83").Value = 3434.4583
$ws.Range("I83").Value = 2954
$ws.Range("K83").Value = 14770
$ws.Range("M83").Value = -9778
# row 97
$ws.Range("H97").Value = 1769.2307
$ws.Range("I97").Value = 1027.7778
$ws.Range("K97").Value = 1027.7778
$ws.Range("M97").Value = -531.7778000000001

$ws = $wb.Worksheets.Item("LTW")
# row 122
$ws.Range("H122").Value = 3939.4902
$ws.Range("I122").Value = 3324.742
$ws.Range("J122").Value = 4892.35
$ws.Range("K122").Value = 9974.226000000001
$ws.Range("L122").Value = 14677.05
$ws.Range("M122").Value = -7524.226000000001
$ws.Range("N122").Value = -19577.05
# row 132
$ws.Range("H132").Value = 3768
$ws.Range("J132").Value = 4072.5557
$ws.Range("L132").Value = 12217.6671
$ws.Range("N132").Value = -17277.6671
# row 136
$ws.Range("H136").Value = 5603.091
$ws.Range("I136").Value = 4493.75
$ws.Range("K136").Value = 13481.25
$ws.Range("M136").Value = -10931.25

$ws = $wb.Worksheets.Item("WVR")
# row 96
$ws.Range("H96").Value = 4200
$ws.Range("I96").Value = 4200
$ws.Range("J96").Value = 4200
$ws.Range("K96").Value = 4200
$ws.Range("L96").Value = 4200
$ws.Range("M96").Value = -2827
$ws.Range("N96").Value = -6946
# row 100
$ws.Range("H100").Value = 1363.2222
$ws.Range("I100").Value = 1382.9524
$ws.Range("K100").Value = 2765.9048
$ws.Range("M100").Value = -2224.9048
# row 132
$ws.Range("H132").Value = 1038.0754
$ws.Range("I132").Value = 891.42224
$ws.Range("K132").Value = 2674.26672
$ws.Range("M132").Value = -144.2667200000001
# row 136
$ws.Range("H136").Value = 8427.875
$ws.Range("I136").Value = 9949.083000000001
$ws.Range("K136").Value = 29847.249
$ws.Range("M136").Value = -27297.249

